$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Formula = "V1.34"
$ws.Range("B6").Formula = "'- No colour sheme / boring colours`n- No dividers between buttons`n- Make errors display in a nicer way"

$ws.Rows(6).RowHeight = 90

$ws.Range("B6").Select()
